$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 27-28 (everything from old row 27 downward shifts down by 2).
$ws.Range("A27:A28").EntireRow.Insert()

# New row 27: Maracuyá, Especial, 20, 60000, Perú
$ws.Range("A27").Value = 9
$ws.Range("B27").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C27").Value = "Metropolitana"
$ws.Range("D27").Value = 44529
$ws.Range("E27").Value = 13
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100108
$ws.Range("H27").Value = "Tropicales y subtropicales"
$ws.Range("I27").Value = 100108003
$ws.Range("J27").Value = "Maracuyá"
$ws.Range("K27").Value = "Sin especificar"
$ws.Range("L27").Value = "Especial"
$ws.Range("M27").Value = 20
$ws.Range("N27").Value = 60000
$ws.Range("O27").Value = 60000
$ws.Range("P27").Value = 60000
$ws.Range("Q27").Value = "`$/caja 18 kilos"
$ws.Range("R27").Value = "Perú"
$ws.Range("S27").Value = 3333
$ws.Range("T27").Value = 18

# New row 28: Maracuyá, Primera, 50, 58000, Perú
$ws.Range("A28").Value = 9
$ws.Range("B28").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C28").Value = "Metropolitana"
$ws.Range("D28").Value = 44529
$ws.Range("E28").Value = 13
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100108
$ws.Range("H28").Value = "Tropicales y subtropicales"
$ws.Range("I28").Value = 100108003
$ws.Range("J28").Value = "Maracuyá"
$ws.Range("K28").Value = "Sin especificar"
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 50
$ws.Range("N28").Value = 58000
$ws.Range("O28").Value = 58000
$ws.Range("P28").Value = 58000
$ws.Range("Q28").Value = "`$/caja 18 kilos"
$ws.Range("R28").Value = "Perú"
$ws.Range("S28").Value = 3222
$ws.Range("T28").Value = 18
